$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column from 2023-09-05 (45174)
# to 2023-09-06 (45175) for rows 2 through 6.
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45175
}
